$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text representation instead of
# being auto-converted to a floating point number by the value-parsing logic.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "65.833.73"
$ws.Range("E2").Value = "  +0.27%  "

$ws.Range("D3").Value = "2.662.63"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "600.03"
$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").Value = "160.27"
$ws.Range("E6").Value = "  +2.40%  "

$ws.Range("D7").Value = "0.642"
$ws.Range("E7").Value = "  +4.24%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -1.97%  "

$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "5.88"
$ws.Range("E10").Value = "  +0.42%  "

$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value = "0.401"
$ws.Range("E11").Value = "  +0.30%  "

$ws.Range("E12").Value = "  +1.61%  "

$ws.Range("D13").Value = "29.23"
$ws.Range("E13").Value = "  -0.25%  "

$ws.Range("E14").Value = "  -0.47%  "

$ws.Range("D15").Value = "3.144.27"
$ws.Range("E15").Value = "  -0.07%  "

$ws.Range("D16").Value = "65.742.45"
$ws.Range("E16").Value = "  +0.37%  "

$ws.Range("D17").Value = "2.607.86"
$ws.Range("E17").Value = "  -1.95%  "

$ws.Range("D18").Value = "12.62"
$ws.Range("E18").Value = "  -1.48%  "

$ws.Range("D19").Value = "4.81"
$ws.Range("E19").Value = "  +0.67%  "

$ws.Range("D20").Value = "355.79"
$ws.Range("E20").Value = "  +1.22%  "

$ws.Range("D21").Value = "7.48"
$ws.Range("E21").Value = "  -1.14%  "

$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("D23").Value = "70.07"
$ws.Range("E23").Value = "  +0.83%  "

$ws.Range("E24").Value = "  +10.02%  "

$ws.Range("D25").Value = "0.0000113"
$ws.Range("E25").Value = "  +1.47%  "

$ws.Range("D26").Value = "9.76"
$ws.Range("E26").Value = "  +1.73%  "

$ws.Range("D27").Value = "1.62"
$ws.Range("E27").Value = "  +2.47%  "

$ws.Range("D28").Value = "578.03"
$ws.Range("E28").Value = "  +9.99%  "

$ws.Range("D29").Value = "8.17"
$ws.Range("E29").Value = "  +2.10%  "

$ws.Range("E30").Value = "  -1.71%  "

$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  -0.02%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "2.16"
$ws.Range("E32").Value = "  +1.31%  "

$ws.Range("E33").Value = "  +4.41%  "

$ws.Range("D34").Value = "6.74"
$ws.Range("E34").Value = "  +4.42%  "

$ws.Range("D35").Value = "5.49"
$ws.Range("E35").Value = "  +0.19%  "

$ws.Range("E36").Value = "  +0.42%  "

$ws.Range("D37").Value = "20.63"
$ws.Range("E37").Value = "  +0.46%  "

$ws.Range("D38").Value = "1.97"
$ws.Range("E38").Value = "  +2.11%  "

$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("D40").Value = "154.22"
$ws.Range("E40").Value = "  -2.49%  "

$ws.Range("D41").Value = "2.48"
$ws.Range("E41").Value = "  +7.32%  "

$ws.Range("D42").Value = "162.20"
$ws.Range("E42").Value = "  -0.65%  "

$ws.Range("D43").Value = "4.11"
$ws.Range("E43").Value = "  -0.16%  "

$ws.Range("D44").Value = "0.0619"
$ws.Range("E44").Value = "  +1.77%  "

$ws.Range("D45").Value = "23.53"
$ws.Range("E45").Value = "  +3.47%  "

$ws.Range("D46").Value = "0.645"
$ws.Range("E46").Value = "  +1.04%  "

$ws.Range("E47").Value = "  +1.05%  "

$ws.Range("E48").Value = "  +1.78%  "

$ws.Range("D49").Value = "19.81"
$ws.Range("E49").Value = "  -1.32%  "

$ws.Range("D50").Value = "0.0₆0247"
$ws.Range("E50").Value = "  -5.99%  "

$ws.Range("D51").Value = "0.820"
$ws.Range("E51").Value = "  +1.36%  "
